# Apply the dataset update: rows 4-13 get new / shifted values, and two new
# rows (11 extra rows total: rows now go to 13) are appended with the most
# recent scraped data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2..13 (columns A..L). Row 1 is the header and is unchanged.
$data = @(
    @(1, 0.4, 1, 1, 1, 0.396, 0.5, 0.896, 200, "Archives", "03-10-2021-15-47", "03-10-2021-15-47"),
    @(0.4285714285714285, 0.5714285714285714, 1, 1, 1, 0.33, 0.5, 0.83, 428, "Archives", "03-11-2021-02-10", "03-11-2021-02-10"),
    @(0.4285714285714285, 0.4285714285714285, 1, 1, 1, 0.3064285714285714, 0.5, 0.8064285714285715, 593, "Archives", "08-11-2021-23-31", "08-11-2021-23-31"),
    @(1, 0.4, 1, 1, 1, 0.396, 0.5, 0.896, 109, "Archives", "13-09-2021-15-00", "13-09-2021-15-00"),
    @(1, 0.4, 1, 1, 1, 0.396, 0.5, 0.896, 112, "Profile", "03-10-2021-15-48", "03-10-2021-15-48"),
    @(0.2, 0.4, 1, 1, 1, 0.264, 0.5, 0.764, 155, "Profile", "08-11-2021-23-31", "08-11-2021-23-31"),
    @(1, 0.4, 1, 1, 1, 0.396, 0.5, 0.896, 112, "Profile", "13-09-2021-14-00", "13-09-2021-14-00"),
    @(0.9090909090909091, 0.04545454545454546, 1, 1, 1, 0.3225, 0.5, 0.8225, 1401, "Frontend", "03-10-2021-15-49", "03-10-2021-15-49"),
    @(0.9591836734693877, 0.02040816326530612, 0.9183673469387755, 1, 0.9387755102040817, 0.3131632653061225, 0.4785714285714285, 0.791734693877551, 4888, "Frontend", "03-11-2021-02-08", "03-11-2021-02-08"),
    @(0.9607843137254902, 0.0196078431372549, 1, 1, 0.9803921568627451, 0.326764705882353, 0.4931372549019607, 0.8199019607843137, 5638, "Frontend", "08-11-2021-23-30", "08-11-2021-23-30"),
    @(0.5, 0.5, 1, 1, 1, 0.33, 0.5, 0.8300000000000001, 35, "Frontend", "13-09-2021-20-00", "13-09-2021-20-00"),
    @(0.9591836734693877, 0.02040816326530612, 0.9183673469387755, 1, 0.9591836734693877, 0.3131632653061225, 0.4857142857142857, 0.7988775510204081, 4858, "Frontend", "17-10-2021-15-30", "17-10-2021-15-30")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
